# The source workbook tracks, per player (column A) and coach/owner
# (column B), whether the player "Started" (column C, Yes/No) on the
# "by Coach" worksheet. This edit flips the Started flag for a specific
# set of rows (the roster's starter designations were corrected/updated).
#
# Target (row -> new "Started" value) is applied explicitly by cell
# address so the result is deterministic regardless of the value
# currently sitting in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

$newStartedValues = @{
    3  = "Yes"
    5  = "No"
    8  = "No"
    9  = "Yes"
    18 = "No"
    20 = "Yes"
    28 = "Yes"
    32 = "No"
    42 = "No"
    45 = "Yes"
    51 = "No"
    54 = "Yes"
    55 = "Yes"
    60 = "No"
    64 = "Yes"
    65 = "Yes"
    67 = "No"
    69 = "No"
    79 = "Yes"
    82 = "No"
}

foreach ($r in $newStartedValues.Keys) {
    $ws.Cells.Item($r, 3).Value2 = $newStartedValues[$r]
}
